# Repull data / push all data / mean calculation
# Updates the "dSF" column (F) values for a set of rows in the single
# worksheet of the workbook to reflect re-pulled source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -5
$ws.Range("F4").Value = 5
$ws.Range("F5").Value = 4
$ws.Range("F7").Value = 4
$ws.Range("F9").Value = -2
$ws.Range("F11").Value = -5
$ws.Range("F12").Value = 1
$ws.Range("F15").Value = -2
$ws.Range("F16").Value = 4
$ws.Range("F17").Value = 6
$ws.Range("F19").Value = -1
$ws.Range("F20").Value = 1
$ws.Range("F21").Value = 1
$ws.Range("F22").Value = 1
$ws.Range("F24").Value = 2
$ws.Range("F25").Value = -2
$ws.Range("F26").Value = 0
$ws.Range("F28").Value = 3
$ws.Range("F29").Value = 2
$ws.Range("F30").Value = -5
$ws.Range("F31").Value = -6
$ws.Range("F32").Value = -1
$ws.Range("F33").Value = 0
$ws.Range("F34").Value = 4
$ws.Range("F36").Value = -3
